$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "243.67"
$ws.Range("D3").Value = "23.29"
$ws.Range("D4").Value = "5.688"
$ws.Range("D5").Value = "0.05819"
$ws.Range("D7").Value = "6.487"
$ws.Range("D9").Value = "0.7993"
$ws.Range("D10").Value = "0.1461"
$ws.Range("D11").Value = "0.07632"
$ws.Range("D12").Value = "0.03262"
$ws.Range("D14").Value = "0.09241"
$ws.Range("D15").Value = "0.001652"
$ws.Range("D16").Value = "3.426"
$ws.Range("D17").Value = "0.04752"
$ws.Range("D18").Value = "0.0006005"
$ws.Range("D19").Value = "0.006229"

# Rows 20 and 21 swap (HotbitToken <-> BitKan), with updated prices
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "0.001072"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "0.003826"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("D23").Value = "3.695"
$ws.Range("D24").Value = "2.209"
$ws.Range("D25").Value = "0.3336"

$ws.Range("D27").Value = "0.0006519"
$ws.Range("E27").Value = "26UpBotsUBXT"

$ws.Range("D40").Value = "0.04323"
$ws.Range("D41").Value = "0.007043"
$ws.Range("D42").Value = "0.1053"

$ws.Range("D44").Value = "0.008619"
$ws.Range("E44").Value = "43LocalTradersLCTWorstin24h"

$ws.Range("D46").Value = "0.00005764"
$ws.Range("D48").Value = "0.7874"

$ws.Range("D49").Value = "0.1048"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"

$ws.Range("D50").Value = "0.00002106"
$ws.Range("D51").Value = "0.01013"
